# Update "Свободное время" schedule sheet:
# - Replace "Арсений" row (row 4) with "Карина" and her own set of time ranges
# - Leave "Андрей" row (row 5) as-is
# - Move active selection to H4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Карина"
$ws.Range("B4").Value = "с 16:00 до 22:00"
$ws.Range("D4").Value = "с 19:00 до 22:00"
$ws.Range("E4").Value = "с 15:00 до 17:00"
$ws.Range("C4").Value = "с 17:00 до 22:00"
$ws.Range("F4").Value = "с 15:00 до 22:00"
$ws.Range("G4").Value = "с 12:00 до 22:00"
$ws.Range("H4").Value = "c 10:00 до 22:00"

$ws.Range("H4").Select()
